# #5: fund, bonds, otherbonds, antique done
#
# Renames sheet5 ("其他有價證券") to "珠寶、古董、字畫" and rebuilds its
# contents as a normalized single-row table (header row + one data row)
# matching the other sheets' "name/quantity/owner/total/property_category/
# category/date/legislator_name/legislator_id/source_file/index" layout.

$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item(5)

# Rename the worksheet tab.
$ws.Name = "珠寶、古董、字畫"

# Wipe the old, malformed contents (A1:G4) entirely before rebuilding.
$ws.Cells.Clear()

# ---- Header row (row 1), columns B:L, bold/centered/bordered like the
# other sheets' header rows (style index 1). ----
$headerCells = @("B1","C1","D1","E1","F1","G1","H1","I1","J1","K1","L1")
$headerValues = @("name","quantity","owner","total","property_category","category","date","legislator_name","legislator_id","source_file","index")

for ($i = 0; $i -lt $headerCells.Length; $i++) {
    $cell = $ws.Range($headerCells[$i])
    $cell.Value = $headerValues[$i]
    $cell.Font.Bold = $true
    $cell.HorizontalAlignment = -4108
    $cell.VerticalAlignment = -4160
    $cell.Borders.LineStyle = 1
}

# ---- Data row (row 2) ----
$ws.Range("A2").Value = 84
$ws.Range("A2").Font.Bold = $true
$ws.Range("A2").HorizontalAlignment = -4108
$ws.Range("A2").VerticalAlignment = -4160
$ws.Range("A2").Borders.LineStyle = 1

$ws.Range("B2").Value = "台鳳高爾夫球場會員證"
$ws.Range("C2").Value = 1
$ws.Range("D2").Value = "蘇震清"
$ws.Range("E2").Value = 760000
$ws.Range("F2").Value = "otherbonds"
$ws.Range("G2").Value = "normal"
$ws.Range("H2").Value = "2011-11-17"
$ws.Range("I2").Value = "蘇震清"
$ws.Range("J2").Value = 1718
$ws.Range("K2").Value = "tmp98701"
$ws.Range("L2").Value = 84
